$d = $word.ActiveDocument

# The template's final conditional block
#   {%p if tool_ID == “divorce_answer” %}
#     ...
#     {%p if there_are_marital_children %}
#       ...
#     {%p else %}
#       ...
#     {%p endif %}
#   (missing closing {%p endif %} for the outer divorce_answer if)
#
# is missing the {%p endif %} that closes the outer
# "if tool_ID == ‘divorce_answer’" block, so the Jinja-style
# if/endif pairs in the template are unbalanced. Locate the empty
# paragraph that should hold that missing endif (it sits between the
# inner endif and the two blank paragraphs that precede the
# "Don’t apply:" heading) and fill it in, matching the formatting
# used by the other "{%p endif %}" markers in the document.

$f = $d.Content.Find
$f.ClearFormatting()
$found = $f.Execute("Don’t apply:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Output "ERROR: anchor text not found"
} else {
    $anchorStart = $f.Parent.Start
    $count = $d.Paragraphs.Count
    $anchorIndex = -1
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $anchorStart -and $p.Range.End -gt $anchorStart) {
            $anchorIndex = $i
        }
    }

    # Two blank paragraphs sit between the inner {%p endif %} and the
    # "Don’t apply:" paragraph; the first of the two blanks is where
    # the missing outer {%p endif %} belongs.
    $target = $d.Paragraphs.Item($anchorIndex - 2)
    $r = $target.Range
    $r.Text = "{%p endif %}"
    $r.Font.NameAscii = "Arial"
    $r.Font.NameOther = "Arial"
    $r.Font.NameBi = "Arial"

    Write-Output "Inserted missing endif at paragraph index $($anchorIndex - 2)"
}
